{"js": "// Office.js (Word JavaScript API) script\n// Applies two edits described by the diff:\n// 1. Merge the \"Monthly Power Generation\" and \" (kWh) \" italic runs into a\n//    single run with text \"Monthly Power Generation (kWh) \".\n// 2. Append two empty (bold-marked) paragraphs followed by a large new\n//    paragraph of body text to the \"sp project\" section, right before the\n//    end of the document body.\n\n// --- Edit 1: merge the two italic runs into one ---------------------------\nconst mergeResults = context.document.body.search(\"Monthly Power Generation (kWh) \", { matchCase: true });\nmergeResults.load(\"items\");\nawait context.sync();\n\nif (mergeResults.items.length > 0) {\n  const mergeRange = mergeResults.items[0];\n  // Re-insert the same text as a straight replace; Word collapses the\n  // underlying run(s) covered by the range into a single run carrying the\n  // formatting of the first run (italic), which merges the two original\n  // runs into one.\n  mergeRange.insertText(\"Monthly Power Generation (kWh) \", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Edit 2: append the new paragraphs at the end of the document --------\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\nconst additionOoxml = `<w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p><w:p><w:r><w:t xml:space=\"preserve\">The models were used to make predictions for </w:t></w:r><w:r><w:t>the</w:t></w:r><w:r><w:t xml:space=\"preserve\"> twelve-month period</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t>at the end of the dataset</w:t></w:r><w:r><w:t xml:space=\"preserve\">. The predictions from the models were compared using </w:t></w:r><w:r><w:t>three performance metrics: RMSE, MAE and R</w:t></w:r><w:r><w:rPr><w:vertAlign w:val=\"superscript\"/></w:rPr><w:t>2</w:t></w:r><w:r><w:t xml:space=\"preserve\">, the results of these are show in Table x. The LSTM-RNN produced the best values for the RMSE and MAE and the </w:t></w:r><w:r><w:t>second-best</w:t></w:r><w:r><w:t xml:space=\"preserve\"> value </w:t></w:r><w:r><w:t>for the R</w:t></w:r><w:r><w:rPr><w:vertAlign w:val=\"superscript\"/></w:rPr><w:t>2</w:t></w:r><w:r><w:t>, with the SARIMAX II model producing the best R</w:t></w:r><w:r><w:rPr><w:vertAlign w:val=\"superscript\"/></w:rPr><w:t>2</w:t></w:r><w:r><w:t xml:space=\"preserve\"> result. </w:t></w:r><w:r><w:t xml:space=\"preserve\">The model predictions </w:t></w:r><w:r><w:t xml:space=\"preserve\">were also compared using visualisations. Figure x shows </w:t></w:r><w:r><w:t xml:space=\"preserve\">each of the predictions made by the models plotted next to the true power output values for that period. </w:t></w:r><w:r><w:t xml:space=\"preserve\">Both SARIMAX models and XGB model without normalisation have a similar pattern to their predictions that mimics the shape of the true values but overestimates the power outputs for the 04/23 </w:t></w:r><w:r><w:t>\u2013</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t>06/23 and the 08</w:t></w:r><w:r><w:t>/23 months. Normalisation appears to stop the XGB model overestimating the power outputs</w:t></w:r><w:r><w:t>. Th</w:t></w:r><w:r><w:t xml:space=\"preserve\">is </w:t></w:r><w:r><w:t>model gives close estimates of power output up to the 06</w:t></w:r><w:r><w:t>/23, then starts to underestimate. The LSTM-RNN model</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t>produces a curve to estimate the power output</w:t></w:r><w:r><w:t>, following the pattern of the true values less closely than the other models, but gives more accurate predictions</w:t></w:r><w:r><w:t xml:space=\"preserve\"> overall. </w:t></w:r><w:r><w:t xml:space=\"preserve\">Figure </w:t></w:r><w:r><w:t>y show</w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t xml:space=\"preserve\"> linear regression for the predicted </w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t xml:space=\"preserve\">onthly </w:t></w:r><w:r><w:t>p</w:t></w:r><w:r><w:t xml:space=\"preserve\">ower </w:t></w:r><w:r><w:t>o</w:t></w:r><w:r><w:t>utputs compare</w:t></w:r><w:r><w:t>d</w:t></w:r><w:r><w:t xml:space=\"preserve\"> against the </w:t></w:r><w:r><w:t>t</w:t></w:r><w:r><w:t xml:space=\"preserve\">rue </w:t></w:r><w:r><w:t>outputs. The R</w:t></w:r><w:r><w:rPr><w:vertAlign w:val=\"superscript\"/></w:rPr><w:t>2</w:t></w:r><w:r><w:t xml:space=\"preserve\"> values, Table x, correspond to these. The more accurate the model predictions, the closer the points lie to the line. Models appear to make better predictions when the true power output</w:t></w:r><w:r><w:t xml:space=\"preserve\"> is low, these low power output values correspond to the winter months.</w:t></w:r><w:r><w:t xml:space=\"preserve\"> However, the winter months for the test set lie</w:t></w:r><w:r><w:t xml:space=\"preserve\"> the closest</w:t></w:r><w:r><w:t xml:space=\"preserve\"> in time to the end of the training set</w:t></w:r><w:r><w:t>, which may be</w:t></w:r><w:r><w:t xml:space=\"preserve\"> causing there more accurate predictions. Future work could use different length horizon times</w:t></w:r><w:r><w:t xml:space=\"preserve\"> for a model</w:t></w:r><w:r><w:t xml:space=\"preserve\"> to evaluate</w:t></w:r><w:r><w:t xml:space=\"preserve\"> whe</w:t></w:r><w:r><w:t>ther</w:t></w:r><w:r><w:t xml:space=\"preserve\"> winter months </w:t></w:r><w:r><w:t xml:space=\"preserve\">or </w:t></w:r><w:r><w:t xml:space=\"preserve\">values closer to the end of the training set are predicted more accurately. </w:t></w:r></w:p>`;\n\nfunction wrapOoxml(bodyFragment) {\n  return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n    '<pkg:xmlData>' +\n    '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n    '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n    '</Relationships>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + bodyFragment + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>';\n}\n\nlastParagraph.insertOoxml(wrapOoxml(additionOoxml), Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script\n# Applies two edits described by the diff:\n# 1. Merge the \"Monthly Power Generation\" and \" (kWh) \" italic runs into a\n#    single run with text \"Monthly Power Generation (kWh) \".\n# 2. Append two empty (bold-marked) paragraphs followed by a large new\n#    paragraph of body text to the end of the \"sp project\" section.\n\n$d = $word.ActiveDocument\n\n# --- Edit 1: merge the two italic runs into one ---------------------------\n$findRange = $d.Content\n$findRange.Find.ClearFormatting()\n$findRange.Find.Replacement.ClearFormatting()\n$findRange.Find.Text = \"Monthly Power Generation (kWh) \"\n$findRange.Find.Replacement.Text = \"Monthly Power Generation (kWh) \"\n$findRange.Find.Execute(\n    [ref]\"Monthly Power Generation (kWh) \",   # FindText\n    [ref]$false,                              # MatchCase\n    [ref]$false,                              # MatchWholeWord\n    [ref]$false,                              # MatchWildcards\n    [ref]$false,                              # MatchSoundsLike\n    [ref]$false,                              # MatchAllWordForms\n    [ref]$true,                               # Forward\n    [ref]1,                                   # Wrap (wdFindContinue)\n    [ref]$false,                              # Format\n    [ref]\"Monthly Power Generation (kWh) \",   # ReplaceWith\n    [ref]2                                    # Replace (wdReplaceAll)\n) | Out-Null\n\n# --- Edit 2: append the new paragraphs at the end of the document --------\n$lastParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$endRange = $d.Range($lastParagraph.Range.End, $lastParagraph.Range.End)\n\n$additionOoxml = @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">\n    <pkg:xmlData>\n      <Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">\n        <Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>\n      </Relationships>\n    </pkg:xmlData>\n  </pkg:part>\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p><w:p><w:r><w:t xml:space=\"preserve\">The models were used to make predictions for </w:t></w:r><w:r><w:t>the</w:t></w:r><w:r><w:t xml:space=\"preserve\"> twelve-month period</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t>at the end of the dataset</w:t></w:r><w:r><w:t xml:space=\"preserve\">. The predictions from the models were compared using </w:t></w:r><w:r><w:t>three performance metrics: RMSE, MAE and R</w:t></w:r><w:r><w:rPr><w:vertAlign w:val=\"superscript\"/></w:rPr><w:t>2</w:t></w:r><w:r><w:t xml:space=\"preserve\">, the results of these are show in Table x. The LSTM-RNN produced the best values for the RMSE and MAE and the </w:t></w:r><w:r><w:t>second-best</w:t></w:r><w:r><w:t xml:space=\"preserve\"> value </w:t></w:r><w:r><w:t>for the R</w:t></w:r><w:r><w:rPr><w:vertAlign w:val=\"superscript\"/></w:rPr><w:t>2</w:t></w:r><w:r><w:t>, with the SARIMAX II model producing the best R</w:t></w:r><w:r><w:rPr><w:vertAlign w:val=\"superscript\"/></w:rPr><w:t>2</w:t></w:r><w:r><w:t xml:space=\"preserve\"> result. </w:t></w:r><w:r><w:t xml:space=\"preserve\">The model predictions </w:t></w:r><w:r><w:t xml:space=\"preserve\">were also compared using visualisations. Figure x shows </w:t></w:r><w:r><w:t xml:space=\"preserve\">each of the predictions made by the models plotted next to the true power output values for that period. </w:t></w:r><w:r><w:t xml:space=\"preserve\">Both SARIMAX models and XGB model without normalisation have a similar pattern to their predictions that mimics the shape of the true values but overestimates the power outputs for the 04/23 </w:t></w:r><w:r><w:t>\u2013</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t>06/23 and the 08</w:t></w:r><w:r><w:t>/23 months. Normalisation appears to stop the XGB model overestimating the power outputs</w:t></w:r><w:r><w:t>. Th</w:t></w:r><w:r><w:t xml:space=\"preserve\">is </w:t></w:r><w:r><w:t>model gives close estimates of power output up to the 06</w:t></w:r><w:r><w:t>/23, then starts to underestimate. The LSTM-RNN model</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t>produces a curve to estimate the power output</w:t></w:r><w:r><w:t>, following the pattern of the true values less closely than the other models, but gives more accurate predictions</w:t></w:r><w:r><w:t xml:space=\"preserve\"> overall. </w:t></w:r><w:r><w:t xml:space=\"preserve\">Figure </w:t></w:r><w:r><w:t>y show</w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t xml:space=\"preserve\"> linear regression for the predicted </w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t xml:space=\"preserve\">onthly </w:t></w:r><w:r><w:t>p</w:t></w:r><w:r><w:t xml:space=\"preserve\">ower </w:t></w:r><w:r><w:t>o</w:t></w:r><w:r><w:t>utputs compare</w:t></w:r><w:r><w:t>d</w:t></w:r><w:r><w:t xml:space=\"preserve\"> against the </w:t></w:r><w:r><w:t>t</w:t></w:r><w:r><w:t xml:space=\"preserve\">rue </w:t></w:r><w:r><w:t>outputs. The R</w:t></w:r><w:r><w:rPr><w:vertAlign w:val=\"superscript\"/></w:rPr><w:t>2</w:t></w:r><w:r><w:t xml:space=\"preserve\"> values, Table x, correspond to these. The more accurate the model predictions, the closer the points lie to the line. Models appear to make better predictions when the true power output</w:t></w:r><w:r><w:t xml:space=\"preserve\"> is low, these low power output values correspond to the winter months.</w:t></w:r><w:r><w:t xml:space=\"preserve\"> However, the winter months for the test set lie</w:t></w:r><w:r><w:t xml:space=\"preserve\"> the closest</w:t></w:r><w:r><w:t xml:space=\"preserve\"> in time to the end of the training set</w:t></w:r><w:r><w:t>, which may be</w:t></w:r><w:r><w:t xml:space=\"preserve\"> causing there more accurate predictions. Future work could use different length horizon times</w:t></w:r><w:r><w:t xml:space=\"preserve\"> for a model</w:t></w:r><w:r><w:t xml:space=\"preserve\"> to evaluate</w:t></w:r><w:r><w:t xml:space=\"preserve\"> whe</w:t></w:r><w:r><w:t>ther</w:t></w:r><w:r><w:t xml:space=\"preserve\"> winter months </w:t></w:r><w:r><w:t xml:space=\"preserve\">or </w:t></w:r><w:r><w:t xml:space=\"preserve\">values closer to the end of the training set are predicted more accurately. </w:t></w:r></w:p></w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n\n$endRange.InsertXML($additionOoxml) | Out-Null\n"}
